# Generate Report for Handback
#
# Updates the localization-status report to reflect a failed handback
# transform for the 51a8e806-... row:
#   - Status text "Ready for handoff" -> "Handback transform failed"
#     on the Overview sheet (E/F columns) and on the zh-cn / de-de
#     per-locale sheets (Status column).
#   - Populates the "Error Detail" column (P) on the zh-cn / de-de
#     sheets with the handback/handoff file-name mismatch message, and
#     widens that column so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 51a8e806-... file; E/F hold the zh-cn / de-de
# status for that row.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn sheet: row 3 ("51a8e806-...") Status column (C) + Error Detail (P).
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("P3").Value = "Handback file name: qm4zmjpu.yoa is different with handoff file name: 51a8e806-d706-408f-8ff8-1478bbd0f44f.53be977b93291c59ca6df56ff6fd7f6ffd3a0563.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# de-de sheet: row 3 ("51a8e806-...") Status column (C) + Error Detail (P).
$dede.Range("C3").Value = $newStatus
$dede.Range("P3").Value = "Handback file name: qm4zmjpu.yoa is different with handoff file name: 51a8e806-d706-408f-8ff8-1478bbd0f44f.53be977b93291c59ca6df56ff6fd7f6ffd3a0563.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.17
